$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.041303753852844
$ws.Range("B1").Value = 2.301737785339355
$ws.Range("C1").Value = 6.771317481994629
$ws.Range("D1").Value = 2.305720329284668
$ws.Range("E1").Value = 1.312609076499939
